# [IMP] z0bug_odoo: test data for MtO
# Add a new "payment_term_id" column (H) to the sale_order test-data sheet,
# populate the first couple of rows, size the new column, and reset the
# active selection back to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header -----------------------------------------------------
$ws.Range("H1").Value = "payment_term_id"

# Give the header the same "arial/9pt/black" look as the other black header
# cell (G1), built off its own formatting so the run stays self-consistent.
$ws.Range("G2").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Font.Name = "arial"

# --- New column data ---------------------------------------------------
$ws.Range("H2").Value = "z0bug.payment_1"
$ws.Range("H1").Copy()
$ws.Range("H2").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("H3").Value = "z0bug.payment_2"
$ws.Range("A2").Copy()
$ws.Range("H3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H3").Font.Color = 0

# --- Column sizing ----------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 15.1

# --- Reset selection to A1 ---------------------------------------------------
$ws.Range("A1").Select() | Out-Null
$excel.CutCopyMode = $false
